# Append two new diary-entry paragraphs (date line + weather line) at the
# very end of the document body, right before the sectPr, mirroring the
# existing "2022年6月Xx日星期x" / weather-description paragraph pairs.

$d = $word.ActiveDocument

# Insertion point: collapsed range at the very end of the document's
# current last paragraph.
$endOfDoc = $d.Paragraphs.Last.Range.End
$insertionPoint = $d.Range($endOfDoc, $endOfDoc)

# Build the two new paragraphs as raw WordprocessingML so the run/para
# property structure (including which runs do/don't carry the
# w:rFonts[@w:hint='eastAsia'] hint) matches exactly, rather than
# inheriting formatting from whatever run currently sits at the end of
# the document.
$newParagraphsXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="eastAsia"/>
    </w:rPr>
    <w:t>2</w:t>
  </w:r>
  <w:r>
    <w:t>022</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="eastAsia"/>
    </w:rPr>
    <w:t>年6月7日星期二</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:hint="eastAsia"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="eastAsia"/>
    </w:rPr>
    <w:t>晴，今天是高考第一天，上午考语文，下午考数学。</w:t>
  </w:r>
</w:p>
"@

$insertionPoint.InsertXML($newParagraphsXml)
